# input.xlsx test-data refresh:
#  - simplify the header row labels (URL / Folder / Title)
#  - turn the 2nd sample URL row into a "missing folder & title" test case
#  - make both sample URLs real clickable hyperlinks
#  - move the active selection to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "URL" / "Folder" / "Title" (was " Video URL"/"Destination Folder"/" Video Title") ---
$ws.Range("A1").Value = "URL"
$ws.Range("B1").Value = "Folder"
$ws.Range("C1").Value = "Title"

# Row 1 no longer needs the taller custom height - let it size back to the default
$ws.Rows.Item(1).AutoFit()

# --- Row 3: clear the Folder/Title sample values -> new "missing folder and/or title" test case ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# --- Turn the two sample URLs into real hyperlinks, preserving their existing (already
#     underlined/hyperlink-colored) cell formatting: Excel's Hyperlinks.Add re-applies the
#     built-in "Hyperlink" style to the target cell, so stash + restore the original format
#     around the call, same as the usual VBA workaround. ---
$ws.Range("A2").Copy($ws.Range("Z1"))
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.youtube.com/watch?v=okpwkeclMu8")
$ws.Range("Z1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A3").Copy($ws.Range("Z1"))
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.youtube.com/watch?v=WvhYuDvH17I")
$ws.Range("Z1").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false

# --- Selection moves from C3 to A3 ---
$ws.Range("A3").Select() | Out-Null
